# ConvertToExtensionBlock fix: employee names are rendered upper-case and the
# "envelope" glyph prefix is dropped from e-mail addresses (the glyph doesn't
# survive the net9 extension-block conversion yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: Full Name -> UPPERCASE
$ws.Range("B2").Value = "JOHN DOE"
$ws.Range("B3").Value = "JANE SMITH"
$ws.Range("B4").Value = "BOB JOHNSON"
$ws.Range("B5").Value = "ALICE BROWN"

# Column C: Email Address -> drop the leading "✉ " glyph/prefix
$ws.Range("C2").Value = "john@company.com"
$ws.Range("C3").Value = "jane@company.com"
$ws.Range("C4").Value = "bob@company.com"
$ws.Range("C5").Value = "alice@company.com"

# Re-fit the two edited columns to their new content width.
$ws.Columns.Item(2).ColumnWidth = 13.33
$ws.Columns.Item(3).ColumnWidth = 19.15
